# Add a "parameters" row-type example (rows=xx) to the oversight template,
# right after the existing blank spacer row that follows the "filter" block
# (old row 50) and before the "action" block (old row 51).
#
# This mirrors inserting three new rows at row 51, which pushes every
# subsequent row down by three (old row 51 -> new row 54, ..., old row 78 ->
# new row 81), and then filling in the three new rows with the new
# "parameters" / "rows=xx" content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert three blank rows before the current row 51 (inserting repeatedly at
# the same index pushes the previously inserted rows down, ending up with
# three blank rows at 51-53, each cloning the formatting of the row above).
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(51).Insert()

# Fill in the new example row describing the "rows=xx" parameter, in the
# same left-to-right, top-to-bottom order the values were originally
# authored so shared-string indices line up the same way.
$ws.Range("C52").Value = "rows=xx"
$ws.Range("D52").Value = "Number of lines to show when editing a text question, xx shuld be replaced by the number"
$ws.Range("B51").Value = "parameters"

# Restore the view: scroll position near the top of the new "parameters"
# rows, with B52 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B52").Select()
